# Applies the cryptos price-list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.254.88"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "2.525.48"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.45"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").Value = "2.524.22"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").Value = "  +1.59%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("E13").Value = "  -3.30%  "

$ws.Range("D14").Value = "2.954.14"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.27"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.83%  "

$ws.Range("D16").Value = "58.956.21"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "2.508.38"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").Value = "0.0₃0784"
$ws.Range("E29").Value = "  +1.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.78"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.93%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.26%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.91%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.829"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.56%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "290.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.12%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.55%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.30"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.606"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "129.53"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.77%  "

$ws.Range("E47").Value = "  +0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0934"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0223"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.51"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.56%  "
